$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before the old "Outstanding" column (N),
# pushing Outstanding/heading/Disbursement one column to the right.
$ws.Columns("N:N").Insert()

# Give the newly inserted column the same width as its left neighbour
# (column M / "In Advance"), matching the look of the rest of the table.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and restore its selection.
$ws.Activate()
$ws.Range("H17").Select()
